# Insert a new data row at row 183 (pushes the existing rows 183-194 down to 184-195)
# and populate it with the new "Choclo" observation for Región del Maule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(183).Insert()

$ws.Range("A183").Value2 = 5
$ws.Range("B183").Value2 = "Macroferia Regional de Talca"
$ws.Range("C183").Value2 = "Maule"
$ws.Range("D183").Value2 = 44610
$ws.Range("E183").Value2 = 7
$ws.Range("F183").Value2 = 100112024
$ws.Range("G183").Value2 = "Choclo"
$ws.Range("H183").Value2 = "Choclero"
$ws.Range("I183").Value2 = "Primera"
$ws.Range("J183").Value2 = 60000
$ws.Range("K183").Value2 = 100
$ws.Range("L183").Value2 = 120
$ws.Range("M183").Value2 = 110
$ws.Range("N183").Value2 = "$/unidad"
$ws.Range("O183").Value2 = "Región del Maule"
$ws.Range("P183").Value2 = 110
$ws.Range("Q183").Value2 = 1
$ws.Range("R183").Value2 = "Hortaliza"
